$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "New Name" -> "New_Name" (A1)
$ws.Range("A1").Value = "New_Name"

# Give the last-row cells in columns A and C the same (non-default) cell
# style already used by the rest of the table instead of the default style.
$ws.Range("A4").Style = $ws.Range("A2").Style
$ws.Range("C4").Style = $ws.Range("A2").Style

# Move the active selection back to A1 (was C4).
$ws.Range("A1").Select()
